$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Fixed log likelihood test to match data now nans removed"
# Rows 23-31 in columns A:C previously held hard-coded literal results of the
# log-likelihood calculation. Re-derive them as live formulas built from the
# per-row constant (the old literal minus the prior/posterior value) plus the
# corresponding prior/posterior cell in rows 1-9, so the sheet recomputes
# correctly now that NaNs have been stripped out of the source data.
$ws.Range("A23").Formula = "=5491+A1"
$ws.Range("B23").Formula = "=5493+B1"

$ws.Range("A24").Formula = "=798+A2"
$ws.Range("B24").Formula = "=1000+B2"
$ws.Range("C24").Formula = "=3693+C2"

$ws.Range("A25").Formula = "=799+A3"
$ws.Range("B25").Formula = "=998+B3"
$ws.Range("C25").Formula = "=3696+C3"

$ws.Range("A26").Formula = "=399+A4"
$ws.Range("B26").Formula = "=399+B4"

$ws.Range("A27").Formula = "=500+A5"
$ws.Range("B27").Formula = "=500+B5"

$ws.Range("A28").Formula = "=1846+A6"
$ws.Range("B28").Formula = "=1847+B6"

$ws.Range("A29").Formula = "=400+A7"
$ws.Range("B29").Formula = "=399+B7"

$ws.Range("A30").Formula = "=500+A8"
$ws.Range("B30").Formula = "=498+B8"

$ws.Range("A31").Formula = "=1849+A9"
$ws.Range("B31").Formula = "=1847+B9"

# Reflect where the user left the view/selection after making the edit:
# scrolled down so row 22 is the top visible row, with A32 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("A32").Select()
